# Update "想去人数" (want-to-go count, column F) figures across the
# 展览 / 演出 / 全部类型 sheets, per the refreshed data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 137
$ws.Range("F3").Value = 386
$ws.Range("F4").Value = 198
$ws.Range("F5").Value = 29
$ws.Range("F7").Value = 453
$ws.Range("F9").Value = 196
$ws.Range("F10").Value = 154
$ws.Range("F11").Value = 179
$ws.Range("F16").Value = 1514
$ws.Range("F17").Value = 552
$ws.Range("F18").Value = 233
$ws.Range("F21").Value = 827
$ws.Range("F22").Value = 1155
$ws.Range("F25").Value = 2665
$ws.Range("F26").Value = 1452
$ws.Range("F27").Value = 66
$ws.Range("F29").Value = 420
$ws.Range("F30").Value = 429
$ws.Range("F31").Value = 1257
$ws.Range("F32").Value = 828
$ws.Range("F33").Value = 1371
$ws.Range("F34").Value = 165
$ws.Range("F36").Value = 788
$ws.Range("F37").Value = 621
$ws.Range("F38").Value = 677
$ws.Range("F39").Value = 853
$ws.Range("F40").Value = 366
$ws.Range("F41").Value = 252

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 633
$ws.Range("F18").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 137
$ws.Range("F6").Value = 386
$ws.Range("F7").Value = 198
$ws.Range("F8").Value = 29
$ws.Range("F12").Value = 453
$ws.Range("F14").Value = 196
$ws.Range("F16").Value = 154
$ws.Range("F17").Value = 179
$ws.Range("F21").Value = 1514
$ws.Range("F22").Value = 552
$ws.Range("F23").Value = 233
$ws.Range("F28").Value = 1155
$ws.Range("F29").Value = 2665
$ws.Range("F30").Value = 1452
$ws.Range("F31").Value = 66
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = 420
$ws.Range("F35").Value = 429
$ws.Range("F36").Value = 1257
$ws.Range("F39").Value = 828
$ws.Range("F40").Value = 1371
$ws.Range("F41").Value = 788
$ws.Range("F42").Value = 621
$ws.Range("F43").Value = 677
$ws.Range("F44").Value = 853
$ws.Range("F45").Value = 366
$ws.Range("F48").Value = 252
